$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows that need the taller (39pt) auto-wrapped height, matching the target
$tallRows = @(303, 310)

$srcFormat = $ws.Range("A300:C300")

# Row 301: lab.mixture.menu
$srcFormat.Copy($ws.Range("A301:C301")) | Out-Null
$ws.Cells.Item(301, 1).Value = "cs"
$ws.Cells.Item(301, 2).Value = 'lab.mixture.menu'
$ws.Cells.Item(301, 3).Value = 'Mix'

# Row 302: lab.mixture.title
$srcFormat.Copy($ws.Range("A302:C302")) | Out-Null
$ws.Cells.Item(302, 1).Value = "cs"
$ws.Cells.Item(302, 2).Value = 'lab.mixture.title'
$ws.Cells.Item(302, 3).Value = 'Mixy'

# Row 303: lab.mixture.subtitle
$srcFormat.Copy($ws.Range("A303:C303")) | Out-Null
$ws.Cells.Item(303, 1).Value = "cs"
$ws.Cells.Item(303, 2).Value = 'lab.mixture.subtitle'
$ws.Cells.Item(303, 3).Value = 'Každý požitek potřebuje liquid a tato sekce slouží pro správu namíchaných liquidů (včetně hotových); tyto mixy se pak dále používají ve vapování pro trasování, jak který mix chutnal.'

# Row 304: lab.mixture.button.create
$srcFormat.Copy($ws.Range("A304:C304")) | Out-Null
$ws.Cells.Item(304, 1).Value = "cs"
$ws.Cells.Item(304, 2).Value = 'lab.mixture.button.create'
$ws.Cells.Item(304, 3).Value = 'Nový mix'

# Row 305: lab.mixture.button.list
$srcFormat.Copy($ws.Range("A305:C305")) | Out-Null
$ws.Cells.Item(305, 1).Value = "cs"
$ws.Cells.Item(305, 2).Value = 'lab.mixture.button.list'
$ws.Cells.Item(305, 3).Value = 'Seznam mixů'

# Row 306: lab.mixture.create.title
$srcFormat.Copy($ws.Range("A306:C306")) | Out-Null
$ws.Cells.Item(306, 1).Value = "cs"
$ws.Cells.Item(306, 2).Value = 'lab.mixture.create.title'
$ws.Cells.Item(306, 3).Value = 'Nový mix'

# Row 307: lab.mixture.create.subtitle
$srcFormat.Copy($ws.Range("A307:C307")) | Out-Null
$ws.Cells.Item(307, 1).Value = "cs"
$ws.Cells.Item(307, 2).Value = 'lab.mixture.create.subtitle'
$ws.Cells.Item(307, 3).Value = 'Mix vám pomůže sledovat složení liquidu, množství nikotinu a jeho stáří (tzn. zrání).'

# Row 308: lab.mixture.list.title
$srcFormat.Copy($ws.Range("A308:C308")) | Out-Null
$ws.Cells.Item(308, 1).Value = "cs"
$ws.Cells.Item(308, 2).Value = 'lab.mixture.list.title'
$ws.Cells.Item(308, 3).Value = 'Seznam mixů'

# Row 309: lab.vape.title
$srcFormat.Copy($ws.Range("A309:C309")) | Out-Null
$ws.Cells.Item(309, 1).Value = "cs"
$ws.Cells.Item(309, 2).Value = 'lab.vape.title'
$ws.Cells.Item(309, 3).Value = 'Vape'

# Row 310: lab.vape.subtitle
$srcFormat.Copy($ws.Range("A310:C310")) | Out-Null
$ws.Cells.Item(310, 1).Value = "cs"
$ws.Cells.Item(310, 2).Value = 'lab.vape.subtitle'
$ws.Cells.Item(310, 3).Value = 'Toto je nejzajímavější část aplikace, kde si můžete trasovat zážitky z vapování, které vám tak umožní zjistit, jaké spirálky, nastavení vzduchu, liquidy (a jejich stáří) vám vyhovují nejvíce.'

# Row 311: lab.vape.button.create
$srcFormat.Copy($ws.Range("A311:C311")) | Out-Null
$ws.Cells.Item(311, 1).Value = "cs"
$ws.Cells.Item(311, 2).Value = 'lab.vape.button.create'
$ws.Cells.Item(311, 3).Value = 'Nový vape'

# Row 312: lab.vape.button.list
$srcFormat.Copy($ws.Range("A312:C312")) | Out-Null
$ws.Cells.Item(312, 1).Value = "cs"
$ws.Cells.Item(312, 2).Value = 'lab.vape.button.list'
$ws.Cells.Item(312, 3).Value = 'Seznam vapů'

# Row 313: lab.vape.create.title
$srcFormat.Copy($ws.Range("A313:C313")) | Out-Null
$ws.Cells.Item(313, 1).Value = "cs"
$ws.Cells.Item(313, 2).Value = 'lab.vape.create.title'
$ws.Cells.Item(313, 3).Value = 'Nový vape'

# Row 314: lab.vape.create.subtitle
$srcFormat.Copy($ws.Range("A314:C314")) | Out-Null
$ws.Cells.Item(314, 1).Value = "cs"
$ws.Cells.Item(314, 2).Value = 'lab.vape.create.subtitle'
$ws.Cells.Item(314, 3).Value = 'Vape je základní stavební kámen pro záznam chutě a požitku z vapování.'

# Row 315: lab.vape.list.title
$srcFormat.Copy($ws.Range("A315:C315")) | Out-Null
$ws.Cells.Item(315, 1).Value = "cs"
$ws.Cells.Item(315, 2).Value = 'lab.vape.list.title'
$ws.Cells.Item(315, 3).Value = 'Seznam vapů'

foreach ($r in $tallRows) {
    $ws.Rows.Item($r).RowHeight = 39
}

$ws.Range("B310").Select() | Out-Null

Write-Host "Done adding rows 301-315"